$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 184.63637
$ws.Range("I6").Value = 133.2
$ws.Range("K6").Value = 399.6
$ws.Range("M6").Value = -287.6
$ws.Range("H32").Value = 3550.111
$ws.Range("I32").Value = 3640.75
$ws.Range("J32").Value = 3477.6
$ws.Range("K32").Value = 3640.75
$ws.Range("L32").Value = 3477.6
$ws.Range("M32").Value = -3314.75
$ws.Range("N32").Value = -4129.6
$ws.Range("H33").Value = 5976.1665
$ws.Range("I33").Value = 7411.643
$ws.Range("K33").Value = 7411.643
$ws.Range("M33").Value = -7182.643
$ws.Range("H86").Value = 69588.87
$ws.Range("I86").Value = 86536.086
$ws.Range("K86").Value = 86536.086
$ws.Range("M86").Value = -85413.086
$ws.Range("H89").Value = 69588.87
$ws.Range("I89").Value = 86536.086
$ws.Range("K89").Value = 432680.43
$ws.Range("M89").Value = -427064.43
$ws.Range("H135").Value = 3018.2144
$ws.Range("I135").Value = 1802.1818
$ws.Range("K135").Value = 16219.6362
$ws.Range("M135").Value = -13684.6362
$ws.Range("H138").Value = 3453.5918
$ws.Range("J138").Value = 3650.2563
$ws.Range("L138").Value = 10950.7689
$ws.Range("N138").Value = -21230.7689

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3060.6667
$ws.Range("I45").Value = 1165
$ws.Range("K45").Value = 1165
$ws.Range("M45").Value = -788
$ws.Range("H61").Value = 3196.889
$ws.Range("I61").Value = 3254
$ws.Range("J61").Value = 2997
$ws.Range("K61").Value = 3254
$ws.Range("L61").Value = 2997
$ws.Range("M61").Value = -3042
$ws.Range("N61").Value = -3421
$ws.Range("H63").Value = 3394.2222
$ws.Range("I63").Value = 2942.6667
$ws.Range("K63").Value = 2942.6667
$ws.Range("M63").Value = -2256.6667
$ws.Range("H66").Value = 3394.2222
$ws.Range("I66").Value = 2942.6667
$ws.Range("K66").Value = 14713.3335
$ws.Range("M66").Value = -11281.3335
$ws.Range("H74").Value = 143681.86
$ws.Range("I74").Value = 143681.86
$ws.Range("K74").Value = 143681.86
$ws.Range("M74").Value = -142807.86
$ws.Range("H77").Value = 143681.86
$ws.Range("I77").Value = 143681.86
$ws.Range("K77").Value = 718409.2999999999
$ws.Range("M77").Value = -714041.2999999999
$ws.Range("H110").Value = 3399
$ws.Range("I110").Value = 3424.5
$ws.Range("J110").Value = 3373.5
$ws.Range("K110").Value = 3424.5
$ws.Range("L110").Value = 3373.5
$ws.Range("M110").Value = -1379.5
$ws.Range("N110").Value = -7463.5
$ws.Range("H136").Value = 3196.889
$ws.Range("I136").Value = 3254
$ws.Range("J136").Value = 2997
$ws.Range("K136").Value = 9762
$ws.Range("L136").Value = 8991
$ws.Range("M136").Value = -7212
$ws.Range("N136").Value = -14091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2566.4443
$ws.Range("I86").Value = 1566.3334
$ws.Range("K86").Value = 1566.3334
$ws.Range("M86").Value = -443.3334
$ws.Range("H89").Value = 2566.4443
$ws.Range("I89").Value = 1566.3334
$ws.Range("K89").Value = 7831.666999999999
$ws.Range("M89").Value = -2215.666999999999
$ws.Range("H94").Value = 9989.25
$ws.Range("I94").Value = 10559.143
$ws.Range("K94").Value = 10559.143
$ws.Range("M94").Value = -10108.143
$ws.Range("H105").Value = 4380.24
$ws.Range("I105").Value = 4231.8423
$ws.Range("K105").Value = 4231.8423
$ws.Range("M105").Value = -2484.8423
$ws.Range("H107").Value = 1301
$ws.Range("I107").Value = 1301
$ws.Range("K107").Value = 1301
$ws.Range("M107").Value = 619
$ws.Range("H134").Value = 2116.0625
$ws.Range("I134").Value = 1902.1578
$ws.Range("J134").Value = 2928.9
$ws.Range("K134").Value = 5706.4734
$ws.Range("L134").Value = 8786.700000000001
$ws.Range("M134").Value = -3171.4734
$ws.Range("N134").Value = -13856.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4038.12
$ws.Range("J31").Value = 8915.166999999999
$ws.Range("L31").Value = 8915.166999999999
$ws.Range("N31").Value = -9505.166999999999
$ws.Range("H34").Value = 4038.12
$ws.Range("J34").Value = 8915.166999999999
$ws.Range("L34").Value = 8915.166999999999
$ws.Range("N34").Value = -9319.166999999999
$ws.Range("H105").Value = 1897.0769
$ws.Range("I105").Value = 1716.4
$ws.Range("J105").Value = 2499.3333
$ws.Range("K105").Value = 1716.4
$ws.Range("L105").Value = 2499.3333
$ws.Range("M105").Value = 30.59999999999991
$ws.Range("N105").Value = -5993.3333
$ws.Range("H107").Value = 403.33334
$ws.Range("I107").Value = 273.33334
$ws.Range("K107").Value = 273.33334
$ws.Range("M107").Value = 1646.66666
$ws.Range("H122").Value = 3601.8333
$ws.Range("I122").Value = 3365.6667
$ws.Range("J122").Value = 3838
$ws.Range("K122").Value = 10097.0001
$ws.Range("L122").Value = 11514
$ws.Range("M122").Value = -7647.000100000001
$ws.Range("N122").Value = -16414

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 450.94446
$ws.Range("I92").Value = 441.81818
$ws.Range("J92").Value = 465.2857
$ws.Range("K92").Value = 1325.45454
$ws.Range("L92").Value = 1395.8571
$ws.Range("M92").Value = -77.45453999999995
$ws.Range("N92").Value = -3891.8571
$ws.Range("H114").Value = 1413.909
$ws.Range("I114").Value = 340.6
$ws.Range("J114").Value = 2308.3333
$ws.Range("K114").Value = 1021.8
$ws.Range("L114").Value = 6924.999899999999
$ws.Range("M114").Value = 2232.2
$ws.Range("N114").Value = -13432.9999
$ws.Range("H129").Value = 609068.1
$ws.Range("I129").Value = 1066.5
$ws.Range("J129").Value = 1065069.4
$ws.Range("K129").Value = 3199.5
$ws.Range("L129").Value = 3195208.2
$ws.Range("M129").Value = 1800.5
$ws.Range("N129").Value = -3205208.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1391.5555
$ws.Range("I97").Value = 1347.1428
$ws.Range("K97").Value = 1347.1428
$ws.Range("M97").Value = -851.1428000000001
$ws.Range("H113").Value = 2693
$ws.Range("I113").Value = 2399
$ws.Range("J113").Value = 2840
$ws.Range("K113").Value = 2399
$ws.Range("L113").Value = 2840
$ws.Range("M113").Value = -229
$ws.Range("N113").Value = -7180
$ws.Range("H122").Value = 1482.1818
$ws.Range("I122").Value = 967.2222
$ws.Range("J122").Value = 3799.5
$ws.Range("K122").Value = 2901.6666
$ws.Range("L122").Value = 11398.5
$ws.Range("M122").Value = -451.6666
$ws.Range("N122").Value = -16298.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3829.0833
$ws.Range("I40").Value = 3696.75
$ws.Range("J40").Value = 3895.25
$ws.Range("K40").Value = 3696.75
$ws.Range("L40").Value = 3895.25
$ws.Range("M40").Value = -3560.75
$ws.Range("N40").Value = -4167.25
$ws.Range("H46").Value = 19484.666
$ws.Range("I46").Value = 24994.154
$ws.Range("J46").Value = 5160
$ws.Range("K46").Value = 24994.154
$ws.Range("L46").Value = 5160
$ws.Range("M46").Value = -24806.154
$ws.Range("N46").Value = -5536
$ws.Range("H122").Value = 4145.108
$ws.Range("I122").Value = 3493.8125
$ws.Range("J122").Value = 4641.3335
$ws.Range("K122").Value = 10481.4375
$ws.Range("L122").Value = 13924.0005
$ws.Range("M122").Value = -8031.4375
$ws.Range("N122").Value = -18824.0005
$ws.Range("H132").Value = 47620
$ws.Range("I132").Value = 66017.89999999999
$ws.Range("K132").Value = 198053.7
$ws.Range("M132").Value = -195523.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652
$ws.Range("H107").Value = 1125.125
$ws.Range("I107").Value = 1066.8334
$ws.Range("K107").Value = 3200.5002
$ws.Range("M107").Value = -1280.5002
$ws.Range("H122").Value = 3344.125
$ws.Range("I122").Value = 3475.5
$ws.Range("K122").Value = 10426.5
$ws.Range("M122").Value = -7976.5
$ws.Range("H136").Value = 4172.067
$ws.Range("I136").Value = 3963.818
$ws.Range("J136").Value = 4744.75
$ws.Range("K136").Value = 11891.454
$ws.Range("L136").Value = 14234.25
$ws.Range("M136").Value = -9341.454000000002
$ws.Range("N136").Value = -19334.25
